$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 85, shifting existing rows 85-212 down to 86-213.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(85, 1).Value = 3
$ws.Cells.Item(85, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44571
$ws.Cells.Item(85, 5).Value = 5
$ws.Cells.Item(85, 6).Value = 100112001
$ws.Cells.Item(85, 7).Value = "Berenjena"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 90
$ws.Cells.Item(85, 11).Value = 7500
$ws.Cells.Item(85, 12).Value = 8000
$ws.Cells.Item(85, 13).Value = 7778
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 16).Value = 130
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
